$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

$ws.Range("A$row").Value = 131261782
$ws.Range("B$row").Value = 99014
$ws.Range("D$row").Value = "VU"
$ws.Range("E$row").Value = 220787
$ws.Range("F$row").Value = "Knärot"
$ws.Range("G$row").Value = "Goodyera repens"
$ws.Range("H$row").Value = "(L.) R. Br."
$ws.Range("I$row").Value = "1"
$ws.Range("J$row").Value = "stjälkar/strån/skott"
$ws.Range("P$row").Value = "Mattsarve, Gtl"
$ws.Range("Q$row").Value = 728130
$ws.Range("R$row").Value = 6370652
$ws.Range("S$row").Value = 10
$ws.Range("T$row").Value = "Gotland"
$ws.Range("U$row").Value = "Gotland"
$ws.Range("V$row").Value = "Gotland"
$ws.Range("W$row").Value = "Gammelgarn"
$ws.Range("X$row").Value = "I-Got-3625"
$ws.Range("Y$row").Value = "2023-10-04"
$ws.Range("AA$row").Value = "2023-10-04"
$ws.Range("AD$row").Value = $false
$ws.Range("AE$row").Value = $false
$ws.Range("AG$row").Value = $false
$ws.Range("AT$row").Value = ""
$ws.Range("AW$row").Value = "Sofia Lund"
$ws.Range("AX$row").Value = "Ingrid Thomasson"
$ws.Range("AY$row").Value = "Floraväkteri Sverige"
